$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$problemText = @"
Data drives so many policy decisions and solutions. Getting access to data can be tricky. This means that there are gaps in the data, especially data used to measure performance of investments over time.

Currently, the only option is to get one-off access from one source. This is expensive and time consuming for everyone.

There are very few trusted, neutral and impartial parties that can bring open government and private data together.

When providing data to government, there is a risk that the government agency receiving the data can react with legislative and regulatory changes impacting the data providers heavily.
"@
$ws.Range("A4").Value = $problemText

$solutionText = @"
A commercially focused Stats NZ business unit allows us to: 
* give value back to the data suppliers, 
* bring competing data sets together for the first time, 
* be trusted to look after the privacy of the data sets, protecting all parties involved.

Data suppliers will provide us with aggregated data based on the use case. We will standardise, confidentialised, and provide assurances to it and then make it accessible according to the various maturity levels of the customers (it could be a data stream, or it could be a basic web application). Customers can provide feedback to help us improve.

We will bring together open government data and private commercial data to create new and improved value.
"@
$ws.Range("B4").Value = $solutionText

$customerSegmentsText = @"
Central Government

Local Government

Crown Entities

Iwi

Council-controlled organisations

Council-controlled trading organisations

Other customer groups discovered through the pilot
"@
$ws.Range("E4").Value = $customerSegmentsText

$costStructureText = @"
The core costs of running of Data Ventures will come from the core team and some business as usual (BAU) costs (such as software subscriptions, stationary, toilet paper, etc.) and hiring a team.

We have levers in the business model which attributes to the running costs of the data brokerage model:

* Short term return to data providers, a fixed amount (either annual or one off) that we provide to offset their initial investment on the work required to shape/process the data sets we require, i.e. to create certain anonymising and classifications based on the data

* Long term return to data providers, a % share of revenue of the use case, access to the aggregated dataset, leverage of DV market research as a value add channel

Marketing and PR are essential to managing perceptions and ensuring stakeholders are appropriately kept up to date. If we get this wrong,  people assume we are doing bad things (which we aren't) e.g. selling Stats NZ data.

Technology costs (marketing / hosting / tools and apps).              
"@
$ws.Range("D6").Value = $costStructureText

$ws.Range("A7:E7").Merge()
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$scoringText = @"
Scoring
Effort: 6 - Enterprise sales in a condensed period of time and the funding cycle associated with government. Limitation in our complexity due to limiting our sales channel to govt only
Acquisition:  6 - We have developed an operating model that considers what's in it for the data providers, DV and the public benefits. It's a case of using this now to develop any further relationships for new and existing data providers for future products. We will have to do some acquisition of either internal resources or contractors or seconded but either way it's not that hard
Risk: 6 - We have mitigation plans of things we have control of. The risk associated is through lack of control over third parties.
Complexity: 6 - Unknown complexities around long term commercial arrangements between DV and customers and DV and data providers.
Value: 8 - DV views the value to bit high, but customers and data providers need help to understand that value.

"@
$ws.Range("A7").Value = $scoringText
